# Applies the "Add files via upload" commit:
#  - Zub_Gel sheet: split the single "Stk. Artikel" line into two
#    (P_Art_S1 / P_Art_S2), split Menge_S into Menge_S1 / Menge_S2,
#    update the Balkonblenden article list and the total-price formula
#    text, and make this sheet the active / selected one.
#  - Brix_Gel_Stab sheet: no longer the selected sheet (selection moves
#    to Zub_Gel); its own cell contents are unchanged.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Zub_Gel")

# --- Insert both new rows first, while row numbers are still easy to reason
#     about, then fill in the text afterwards. ---

# New second "Stk. Artikel" row, right after the existing one (old row 3)
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).Value = "Mehrfach"
$ws.Cells.Item(4, 2).Value = "Stk. Artikel"

# New second "Menge / Stück" row, right after the (still single) existing
# "Menge / Stück" row, which is now row 6
$ws.Rows.Item(7).Insert()
$ws.Cells.Item(7, 1).Value = "Zahl"

# --- Now update the text content, in the same order the strings were
#     originally authored so the shared-string table layout matches. ---

# Row 2: "Lfm. Artikel" / P_Art -> update the options list (Balkonblenden)
$ws.Cells.Item(2, 4).Value = "BalkonblendenBV 160mm 2 Latten(lfm):33, BalkonblendenBV 240mm 3 Latten(lfm):41, BalkonblendenBV 320mm 4 Latten(lfm):50,BalkonblendenBV 400mm  Latten(lfm):58,"

# Row 3: first "Stk. Artikel" -> options list gains two new Balkonblenden
# accessory entries, variable renamed P_Art_S -> P_Art_S1
$ws.Cells.Item(3, 4).Value = "Blumenkistenhalter (Paar):45, T-Verbinder für Handlauf:25, Eckverbinder 90°:35, Wandanschluss Handlauf:18, Abdeckkappe Steher (Stk):8, Stoßverbinder (Stk):12, Sichtschutz-Element (Stk):250, Seitenblende (Stk):180, Balkonblenden-Halterung für je 1 BV 160-400:26, Balkonblenden-Abwinkelung für BV 160 - 400:60"
$ws.Cells.Item(3, 3).Value = "P_Art_S1"

# Row 4 (new): second "Stk. Artikel" line, variable P_Art_S2, same options
$ws.Cells.Item(4, 3).Value = "P_Art_S2"
$ws.Cells.Item(4, 4).Value = "Blumenkistenhalter (Paar):45, T-Verbinder für Handlauf:25, Eckverbinder 90°:35, Wandanschluss Handlauf:18, Abdeckkappe Steher (Stk):8, Stoßverbinder (Stk):12, Sichtschutz-Element (Stk):250, Seitenblende (Stk):180, Balkonblenden-Halterung für je 1 BV 160-400:26, Balkonblenden-Abwinkelung für BV 160 - 400:60"

# Row 8 (old row 6, "Preis / Gesamtpreis / Endpreis"): update the formula text
$ws.Cells.Item(8, 5).Value = "(P_Art * Menge_L) + (P_Art_S1 * Menge_S1) + (P_Art_S2 * Menge_S2)"

# Row 6 (old row 5, "Zahl" / "Menge / Stück"): rename to the "...1" variant
$ws.Cells.Item(6, 2).Value = "Menge / Stück1"

# Row 7 (new): second "Menge / Stück" row
$ws.Cells.Item(7, 2).Value = "Menge / Stück2"

# Set the variable names for the two quantity rows (order matters for the
# shared-string table layout: Menge_S2 is written before Menge_S1)
$ws.Cells.Item(7, 3).Value = "Menge_S2"
$ws.Cells.Item(6, 3).Value = "Menge_S1"

# Column widths roughly matching the new content (cosmetic best-fit)
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).ColumnWidth = 23.5703125
$ws.Columns.Item(5).AutoFit()

# Make "Zub_Gel" the active sheet / selected cell, matching the new
# activeTab + tabSelected state (moves off "Brix_Gel_Stab")
$ws.Activate()
$ws.Range("D8").Select()
